$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.509.64'
$ws.Range("E2").Value = '  -2.87%  '

# Row 3
$ws.Range("D3").Value = '1.994.94'
$ws.Range("E3").Value = '  -4.70%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.32%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.98'
$ws.Range("E5").Value = '  -4.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.016'
$ws.Range("E6").Value = '  +1.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4979'
$ws.Range("E7").Value = '  -4.89%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4206'
$ws.Range("E8").Value = '  -4.77%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.25'
$ws.Range("E9").Value = '  -0.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08928'
$ws.Range("E10").Value = '  -4.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.107'
$ws.Range("E11").Value = '  -5.20%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.08'
$ws.Range("E12").Value = '  -7.00%  '

# Row 13
$ws.Range("D13").Value = '1.989.79'
$ws.Range("E13").Value = '  -6.07%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.913'
$ws.Range("E14").Value = '  -7.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.406'
$ws.Range("E15").Value = '  -7.16%  '

# Row 16
$ws.Range("E16").Value = '  +1.26%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.74'
$ws.Range("E17").Value = '  -7.51%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001104'
$ws.Range("E18").Value = '  -4.70%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06672'
$ws.Range("E19").Value = '  +0.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.41'
$ws.Range("E20").Value = '  -8.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.010'
$ws.Range("E21").Value = '  +0.91%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.893'
$ws.Range("E22").Value = '  -6.80%  '

# Row 23
$ws.Range("D23").Value = '29.552.44'
$ws.Range("E23").Value = '  -2.80%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.91'
$ws.Range("E24").Value = '  -5.04%  '

# Row 25
$ws.Range("E25").Value = '  -0.22%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.56'
$ws.Range("E26").Value = '  -3.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.61'
$ws.Range("E27").Value = '  -5.53%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.255'
$ws.Range("E28").Value = '  -8.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.281'
$ws.Range("E29").Value = '  -8.86%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.98'
$ws.Range("E30").Value = '  -4.63%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.043'
$ws.Range("E31").Value = '  -8.29%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09887'
$ws.Range("E32").Value = '  -5.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.542'
$ws.Range("E33").Value = '  -7.04%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.809'
$ws.Range("E34").Value = '  -1.37%  '

# Row 35
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.779'
$ws.Range("E35").Value = '  -7.57%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02444'
$ws.Range("E36").Value = '  -7.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.243'
$ws.Range("E37").Value = '  -8.73%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.299'
$ws.Range("E38").Value = '  -3.23%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06362'
$ws.Range("E39").Value = '  -7.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6508'
$ws.Range("E40").Value = '  -6.80%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.54'
$ws.Range("E41").Value = '  -7.97%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2029'
$ws.Range("E42").Value = '  -8.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.013'
$ws.Range("E43").Value = '  +1.21%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6290'
$ws.Range("E44").Value = '  -7.61%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.49'
$ws.Range("E45").Value = '  -6.29%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.176'
$ws.Range("E46").Value = '  -7.17%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.305'
$ws.Range("E47").Value = '  -5.25%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.504'
$ws.Range("E48").Value = '  -3.62%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000339'
$ws.Range("E49").Value = '  -2.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06949'
$ws.Range("E50").Value = '  -4.03%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.122'
$ws.Range("E51").Value = '  -7.90%  '
